$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (38) down into the two
# new rows (39, 40) so the new rows pick up the same column styles/number
# formats without creating new style entries.
$ws.Range("A38:F38").Copy()
$ws.Range("A39:F39").PasteSpecial(-4122)
$ws.Range("B38:F38").Copy()
$ws.Range("B40:F40").PasteSpecial(-4122)

# Row 39: new day 4/2/2020, job 1 - Skype meeting with Hassan, 15:45-16:15, 0.5 hrs
$ws.Range("A39").Value = 43923
$ws.Range("B39").Value = "1"
$ws.Range("C39").Value = "Skype meeting with Hassan"
$ws.Range("D39").Value = 0.65625
$ws.Range("E39").Value = 0.67708333333333337
$ws.Range("F39").Value = 0.5

# Row 40: same day, job 2 - Travel time and risk, 16:20-17:00, 0.5 hrs
$ws.Range("B40").Value = "2"
$ws.Range("C40").Value = "Travel time and risk"
$ws.Range("D40").Value = 0.68055555555555547
$ws.Range("E40").Value = 0.70833333333333337
$ws.Range("F40").Value = 0.5

$ws.Range("F40").Select() | Out-Null
